$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New receptionist row (Constanza) - add this first so the new shared
# string "Constanza" gets allocated before "18_02_2024".
$ws.Range("A6").Value = "Constanza"
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").ClearFormats()
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0

# New week column (18/02/2024)
$ws.Range("I1").Value = "18_02_2024"
$ws.Range("I2").Value = 1164
$ws.Range("I3").Value = 1221
$ws.Range("I4").Value = 1212
$ws.Range("I5").Value = 2989
$ws.Range("I6").Value = 27

$ws.Range("I7").Select() | Out-Null
